$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 184, shifting rows 184:230 down to 185:231
$ws.Rows("184:184").Insert()

# Populate the newly inserted row 184 with the new data record
$ws.Cells.Item(184, 1).Value = 8
$ws.Cells.Item(184, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44551
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 6).Value = 100112032
$ws.Cells.Item(184, 7).Value = "Zapallo italiano"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 480
$ws.Cells.Item(184, 11).Value = 8000
$ws.Cells.Item(184, 12).Value = 9000
$ws.Cells.Item(184, 13).Value = 8500
$ws.Cells.Item(184, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(184, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(184, 16).Value = 121
$ws.Cells.Item(184, 17).Value = 70
$ws.Cells.Item(184, 18).Value = "Hortaliza"
